# feat(master-data): add parsers for master data
#
# Update the input-date header placeholder so that both the submitted
# input date and the (optional) paper input date are rendered through
# the "dd.MM.YYYY" date filter. The three separate template runs that
# made up the old Jinja expression are merged into a single run holding
# the new combined expression.

$d = $word.ActiveDocument

$old = "{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}"
$new = '{% if inputDateHeader %}{{ inputDateHeader | date("dd.MM.YYYY") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date("dd.MM.YYYY") }}){% else %}{% endif %}{% else %}-{% endif %}'

$replaced = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "$old*") {
        $start = $p.Range.Start
        $r = $d.Range($start, $start + $old.Length)
        $r.Text = $new
        $replaced = $replaced + 1
    }
}

Write-Output "Replaced $replaced occurrence(s)."
